# Apply "enemy_strength_increase" (column H) updates to the Locations sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the H value is being removed entirely
$ws.Range("H37").ClearContents()
$ws.Range("H39").ClearContents()

# Rows where an H value is being added or changed
$ws.Range("H40").Value = 0.15
$ws.Range("H41").Value = 0.15
$ws.Range("H42").Value = 0.15
$ws.Range("H44").Value = 0.4
$ws.Range("H45").Value = 0.4
$ws.Range("H46").Value = 0.45
$ws.Range("H47").Value = 0.15
$ws.Range("H48").Value = 0.15
$ws.Range("H49").Value = 0.15
$ws.Range("H63").Value = 0.15
$ws.Range("H64").Value = 0.15
